$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Dhh"
$ws.Range("C2").Value = "Hhip"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.179395
$ws.Range("H2").Value = 6.538185
$ws.Range("I2").Value = 0.4845018986408914
$ws.Range("J2").Value = 0.4845018986408914
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.605175333333333
$ws.Range("N2").Value = 19.815526
$ws.Range("O2").Value = 0.5779586116201393
$ws.Range("P2").Value = 0.5779586116201394
$ws.Range("Q2").Value = 14.39528609559
$ws.Range("R2").Value = 129.55757486031
$ws.Range("S2").Value = 0.280022044665811
$ws.Range("T2").Value = 0.2800220446658111

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Dhh"
$ws.Range("C3").Value = "Hhip"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.179395
$ws.Range("H3").Value = 6.538185
$ws.Range("I3").Value = 0.4845018986408914
$ws.Range("J3").Value = 0.4845018986408914
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.823282
$ws.Range("N3").Value = 14.469846
$ws.Range("O3").Value = 0.4220413883798607
$ws.Range("P3").Value = 0.4220413883798607
$ws.Range("Q3").Value = 10.51183667439
$ws.Range("R3").Value = 94.60653006951001
$ws.Range("S3").Value = 0.2044798539750803
$ws.Range("T3").Value = 0.2044798539750803

# Row 4
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Dhh"
$ws.Range("C4").Value = "Hhip"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7348883333333333
$ws.Range("H4").Value = 2.204665
$ws.Range("I4").Value = 0.1633732264179005
$ws.Range("J4").Value = 0.1633732264179005
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.605175333333333
$ws.Range("N4").Value = 19.815526
$ws.Range("O4").Value = 0.5779586116201393
$ws.Range("P4").Value = 0.5779586116201394
$ws.Range("Q4").Value = 4.854066292087777
$ws.Range("R4").Value = 43.68659662878999
$ws.Range("S4").Value = 0.09442296311639242
$ws.Range("T4").Value = 0.09442296311639244

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Dhh"
$ws.Range("C5").Value = "Hhip"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7348883333333333
$ws.Range("H5").Value = 2.204665
$ws.Range("I5").Value = 0.1633732264179005
$ws.Range("J5").Value = 0.1633732264179005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.823282
$ws.Range("N5").Value = 14.469846
$ws.Range("O5").Value = 0.4220413883798607
$ws.Range("P5").Value = 0.4220413883798607
$ws.Range("Q5").Value = 3.544573670176666
$ws.Range("R5").Value = 31.90116303159
$ws.Range("S5").Value = 0.06895026330150805
$ws.Range("T5").Value = 0.06895026330150805

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Dhh"
$ws.Range("C6").Value = "Hhip"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.583934333333333
$ws.Range("H6").Value = 4.751803
$ws.Range("I6").Value = 0.3521248749412083
$ws.Range("J6").Value = 0.3521248749412082
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.605175333333333
$ws.Range("N6").Value = 19.815526
$ws.Range("O6").Value = 0.5779586116201393
$ws.Range("P6").Value = 0.5779586116201394
$ws.Range("Q6").Value = 10.46216398815311
$ws.Range("R6").Value = 94.15947589337799
$ws.Range("S6").Value = 0.2035136038379359
$ws.Range("T6").Value = 0.2035136038379359

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Dhh"
$ws.Range("C7").Value = "Hhip"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.583934333333333
$ws.Range("H7").Value = 4.751803
$ws.Range("I7").Value = 0.3521248749412083
$ws.Range("J7").Value = 0.3521248749412082
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.823282
$ws.Range("N7").Value = 14.469846
$ws.Range("O7").Value = 0.4220413883798607
$ws.Range("P7").Value = 0.4220413883798607
$ws.Range("Q7").Value = 7.639761959148666
$ws.Range("R7").Value = 68.757857632338
$ws.Range("S7").Value = 0.1486112711032723
$ws.Range("T7").Value = 0.1486112711032723
